$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the previously-blank score cells ---
$ws.Range("I16").Copy()
$ws.Range("J16").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("J16").Value = 5

$ws.Range("I20").Copy()
$ws.Range("J20").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("J20").Value = 5

$ws.Range("G24").Value = 5
$ws.Range("H24").Value = 5

# --- Add the new "total" column L with SUM formulas ---
$ws.Range("L4").Formula = "=SUM(C4:J4)"
$ws.Range("L5:L32").Formula = "=SUM(C5:J5)"

# --- Add color-scale conditional formatting to the new L column ---
$cf = $ws.Range("L4:L32").FormatConditions.AddColorScale(3)
$cf.SetFirstPriority()

# --- Restore the view/selection state ---
$ws.Range("K16").Select()
